$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression (A2 stays the same, B2 stays the same)
$ws.Range("C2").Value = 4906221268062289
$ws.Range("D2").Value = 4906221268062289

# Row 3 - RandomForestRegressor (A3 stays the same)
$ws.Range("B3").Value = 2150667348993121
$ws.Range("C3").Value = 2418470148766574
$ws.Range("D3").Value = 3968113234703352

# Row 4 - GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 5647969298398.135
$ws.Range("C4").Value = 5440060958795.636
$ws.Range("D4").Value = 4380997681772726

# Row 5 - AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 346179985220238.2
$ws.Range("C5").Value = 221583517801467.2
$ws.Range("D5").Value = 566257974029574.8
